$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J header value (row 1)
$ws.Range("J1").Value = 8
$ws.Range("J1").Style = $ws.Range("I1").Style

# Updated / new values for rows 3-7, columns E through J
$data = @{
    3 = @(65.75757575757578, 76.45454545454548, 75.62121212121212, 95.63636363636358, 96.09090909090905, 95.95454545454541)
    4 = @(57.4949494949495, 71.53535353535351, 71.17171717171715, 94.7474747474747, 94.92929292929288, 95.07070707070702)
    5 = @(35.30201342281879, 43.04697986577183, 37.48993288590604, 68.60402684563758, 68.63087248322147, 66.93959731543625)
    6 = @(39.40625, 54.04687500000001, 47.4375, 85.796875, 85.9375, 85.703125)
    7 = @(37.46788990825689, 52.49541284403671, 47.8348623853211, 91.81651376146789, 91.50458715596329, 91.57798165137613)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("E$row").Value = $vals[0]
    $ws.Range("F$row").Value = $vals[1]
    $ws.Range("G$row").Value = $vals[2]
    $ws.Range("H$row").Value = $vals[3]
    $ws.Range("I$row").Value = $vals[4]
    $ws.Range("J$row").Value = $vals[5]
}
